# Apply the "break out stock.yaml completed" edit to the "day" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# 1) D200:D206 currently store bsecode as text (inlineStr). Convert them to
#    real numbers, keeping the same digits.
$bsecodes = @{
    200 = 540115
    201 = 500114
    202 = 500331
    203 = 500325
    204 = 542726
    205 = 500038
    206 = 532523
}
foreach ($row in $bsecodes.Keys) {
    $ws.Cells.Item($row, 4).Value = $bsecodes[$row]
}

# 2) Append the new rows (207-212) reported on 24/07/2024.
$newRows = @(
    @{ A=1; B="HAL";        C="Hindustan Aeronautics Ltd";               D="541154"; E=-0.14; F=4849.5;   G=2601966;  H="day"; I="24/07/2024 11:35:38" },
    @{ A=2; B="KOTAKBANK";  C="Kotak Mahindra Bank Limited";             D="500247"; E=-1.28; F=1746.2;   G=7487207;  H="day"; I="24/07/2024 11:35:38" },
    @{ A=3; B="CHAMBLFERT"; C="Chambal Fertilizers & Chemicals Limited"; D="500085"; E=1.63;  F=494.05;   G=3166226;  H="day"; I="24/07/2024 11:35:38" },
    @{ A=4; B="HINDCOPPER"; C="Hindustan Copper Limited";                D="513599"; E=-0.61; F=310.85;   G=3700107;  H="day"; I="24/07/2024 11:35:38" },
    @{ A=5; B="BHEL";       C="Bharat Heavy Electricals Limited";        D="500103"; E=0.59;  F=309.15;   G=15750887; H="day"; I="24/07/2024 11:35:38" },
    @{ A=6; B="NATIONALUM"; C="National Aluminium Company Limited";      D="532234"; E=0.76;  F=186.57;   G=7595394;  H="day"; I="24/07/2024 11:35:38" }
)

$startRow = 207
$endRow = $startRow + $newRows.Count - 1

# The bsecode column (D) for these new rows stays text (e.g. "541154"), just
# like the rest of column D was before this edit. Force text formatting
# first so the numeric-looking strings aren't auto-converted to numbers.
$ws.Range("D$startRow`:D$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData.A
    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 6).Value = $rowData.F
    $ws.Cells.Item($r, 7).Value = $rowData.G
    $ws.Cells.Item($r, 8).Value = $rowData.H
    $ws.Cells.Item($r, 9).Value = $rowData.I
}
